$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Campos del formulario")

$ws.Range("E1").Value = "Finding"

$comment = $ws.Range("E1").Comment
$comment.Text("Finding (Hallazgo):" + "`nEscribe el numero de la posicion")

$ws.Range("E17:E18").Select()
